# BIS-769: Fixed xls test files
# The DATASET_TYPE export template grows two new columns ("Pattern" and
# "Pattern Type") appended after the existing "Unique" column (L) in the
# property-row header (row 4), i.e. M4 = "Pattern", N4 = "Pattern Type".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing "Unique" header cell (L4) onto the two new header cells
# so they pick up the same bold header style instead of the default style.
$ws.Range("L4").Copy()
$ws.Range("M4:N4").PasteSpecial(-4122)

# Set the new header labels.
$ws.Range("M4").Value = "Pattern"
$ws.Range("N4").Value = "Pattern Type"

# Match the author's final selection/active cell (M4:N4).
$ws.Range("M4:N4").Select() | Out-Null
